$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Good morning. I have an issue with recharge my mobile. When I use
#    recharge card, it shows an error message." ->
#    "Good morning. I have an issue with recharging my account. When I
#    use recharge card, it shows an error message."
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Good morning. I have an issue with recharge my mobile. When I use recharge card, it shows an error message.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Good morning. I have an issue with recharging my account. When I use recharge card, it shows an error message.",
    2)

# ---------------------------------------------------------------------
# 2. Add an empty "_GoBack" bookmark right after "number" (and before
#    the trailing gramEnd proof-error marker) in the paragraph
#    "Could you please tell me your mobile number".
#    The engine mis-resolves a zero-length Range that sits exactly on
#    a paragraph's last character, so we briefly pad the text with two
#    throwaway characters to move that boundary, insert the bookmark at
#    the now-safe offset, then strip the padding back out.
# ---------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("Could you please tell me your mobile number")
$lastChar = $d.Range($anchor.End - 1, $anchor.End)
$lastChar.InsertAfter("ZZ")

$bmSpot = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $bmSpot)

$padding = $d.Content
$padding.Find.Execute("ZZ")
$padding.Text = ""

# ---------------------------------------------------------------------
# 3. Collapse the "Thank you sir. Please leave a feedback after this. "
#    paragraph (currently split across three runs around a gramStart/
#    gramEnd proof-error pair) into a single plain run.
# ---------------------------------------------------------------------
$thank = $d.Content
$thank.Find.Execute("Thank you sir. Please leave a feedback after this. ")
$thank.Text = "__TMP_PLACEHOLDER__"

$thank2 = $d.Content
$thank2.Find.Execute("__TMP_PLACEHOLDER__")
$thank2.Text = "Thank you sir. Please leave a feedback after this. "
